$d = $word.ActiveDocument

# The site footer ("Ver no Jupiter ..." and the copyright notice), along with the blank
# paragraph that separated it from the "Requisitos" list, needs to be removed. It sat right
# after the last requisite line ("LOB1225: Tratamento de Aguas Residuarias (Requisito fraco)")
# and right before the trailing blank paragraph / page-break paragraph that close the document.

$copyrightIdx = -1
$jupiterIdx = -1
$lob1225Idx = -1

for ($i = $d.Paragraphs.Count; $i -ge 1; $i--) {
    $t = $d.Paragraphs.Item($i).Range.Text
    if ($t -like "*Contact: luizeleno@usp.br*" -and $copyrightIdx -eq -1) {
        $copyrightIdx = $i
    }
    if ($t -like "*Ver no Jupiter*" -and $jupiterIdx -eq -1) {
        $jupiterIdx = $i
    }
    if ($t -like "*LOB1225*" -and $lob1225Idx -eq -1) {
        $lob1225Idx = $i
    }
}

# The blank separator paragraph is immediately after the LOB1225 requisite paragraph.
$emptyIdx = $lob1225Idx + 1

# Delete starting from the last paragraph so the indices of the earlier ones stay valid.
if ($copyrightIdx -ne -1) {
    $d.Paragraphs.Item($copyrightIdx).Range.Delete()
}
if ($jupiterIdx -ne -1) {
    $d.Paragraphs.Item($jupiterIdx).Range.Delete()
}
if ($emptyIdx -ne -1) {
    $d.Paragraphs.Item($emptyIdx).Range.Delete()
}
